$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -3
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -3
$ws.Range("F15").Value = 3
$ws.Range("F19").Value = 11
$ws.Range("F22").Value = -5
$ws.Range("F23").Value = 9
$ws.Range("F28").Value = -3
$ws.Range("F29").Value = 0
